$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-19 13:05:48"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-19 13:05:44"
$wsZhCn.Range("K2").Value = "2016-08-19 13:06:06"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-19 13:05:48"
$wsDeDe.Range("K2").Value = "2016-08-19 13:06:21"
